$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.262.19"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.239.19"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "293.62"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "88.35"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.513"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.474"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.37"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0783"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.51"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.584.32"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.86"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.214.11"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.735"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "40.173.60"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0891"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.43"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +6.96%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.70"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.24"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.82"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.87"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.25"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "155.13"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.19"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.71%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.94"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0720"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.36"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.90"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +6.47%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.82"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0976"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.70"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.133.92"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +5.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.86"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "18.35"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +11.14%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.61%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.86"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.69"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.446.41"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.49"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.78%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "89.08"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.77%  "

Write-Host "Updated cryptos list"
